$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 2 (top of data, before "J5, J9") for the RESET entry.
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = "431256083736"
$ws.Cells.Item(2, 2).Value = "RESET"
$ws.Cells.Item(2, 3).Value = "4312560837X6"

# Extend the designator list for the 13K / R1206 resistor row.
$ws.Cells.Item(14, 2).Value = "R14, R15, R16, R17, R19, R20, R23, R26, R31, R32, R33, R34, R41, R44, R53, R56"

# Insert a new row after the 1K2/R21.. row for the additional 1k2/R0805 resistor group.
$ws.Rows.Item(17).Insert()
$ws.Cells.Item(17, 1).Value = "1k2"
$ws.Cells.Item(17, 2).Value = "R35, R36, R37, R38"
$ws.Cells.Item(17, 3).Value = "R0805"

# Insert a new row after the 645004114822/J2 row for the diode/GS1J-L entry.
$ws.Rows.Item(31).Insert()
$ws.Cells.Item(31, 1).Value = "GS1J-L"
$ws.Cells.Item(31, 2).Value = "D014, D016, D018, D019, D020, D021, D022, D023"
$ws.Cells.Item(31, 3).Value = "SMA_DO-214AC"
